## Allow summary columns to reference different tables (#11)
## * Load tables from all sheets and clarify summary sheet flag
## * Removed duplicate columns from linked tables
##
## Adds a second worksheet ("Sheet2") holding a small "headcount" table
## (Table2: Category / Perm Employee / Contract Employee) and wires a new
## "Total Employee" summary column on Sheet1 that pulls from Table2 instead
## of Raw_Data, mirroring the existing FILTER-based summary columns.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. New "Sheet2" with the Perm/Contract headcount data, right after Sheet1
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Category"
$ws2.Range("B1").Value = "Perm Employee"
$ws2.Range("C1").Value = "Contract Employee"

$ws2.Range("A2").Value = "Electronics"
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 6

$ws2.Range("A3").Value = "Home"
$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = 4

$ws2.Range("A4").Value = "Wearables"
$ws2.Range("B4").Value = 1
$ws2.Range("C4").Value = 3

$ws2.Range("A5").Value = "Footwear"
$ws2.Range("B5").Value = 1
$ws2.Range("C5").Value = 2

# Turn the range into a proper table ("Table2") styled like Raw_Data
$tbl2 = $ws2.ListObjects.Add(1, $ws2.Range("A1:C5"), $null, 1)
$tbl2.Name = "Table2"
$tbl2.TableStyle = "TableStyleMedium9"

[void]$ws2.Range("A1:C5").Select()

# ---------------------------------------------------------------------
# 2. New "Total Employee" summary column (G) on Sheet1, driven off Table2
# ---------------------------------------------------------------------
$ws1.Range("G11").Value = "Total Employee"

$ws1.Range("G12").FormulaArray = "=SUM(FILTER(Table2[Perm Employee],Table2[Category]=A12))+SUM(FILTER(Table2[Contract Employee],Table2[Category]=A12))"
$ws1.Range("G13").FormulaArray = "=SUM(FILTER(Table2[Perm Employee],Table2[Category]=A13))+SUM(FILTER(Table2[Contract Employee],Table2[Category]=A13))"
$ws1.Range("G14").FormulaArray = "=SUM(FILTER(Table2[Perm Employee],Table2[Category]=A14))+SUM(FILTER(Table2[Contract Employee],Table2[Category]=A14))"
$ws1.Range("G15").FormulaArray = "=SUM(FILTER(Table2[Perm Employee],Table2[Category]=A15))+SUM(FILTER(Table2[Contract Employee],Table2[Category]=A15))"

# Leftover widened columns (L:N) from prototyping the headcount table on Sheet1
# before it was moved to its own sheet.
$ws1.Columns.Item(12).ColumnWidth = 10.44
$ws1.Columns.Item(13).ColumnWidth = 16.11
$ws1.Columns.Item(14).ColumnWidth = 18.89

# ---------------------------------------------------------------------
# 3. Restore Sheet1 as the active sheet / selection, as in the saved file
# ---------------------------------------------------------------------
$ws1.Activate()
[void]$ws1.Range("B17").Select()
